# Applies the "Added a few more slots" edit:
#  1. Removes the standalone "Meta description: ..." paragraph that used to
#     sit right under the title.
#  2. Adds a new bold "Play Fruit Awards | Free Review of Classic Slot Game"
#     paragraph right before the closing (italic) paragraph, and rewrites
#     that closing paragraph's text to the meta-description sentence
#     (keeping its italic formatting).

$d = $word.ActiveDocument

# --- Step 1: drop the "Meta description" paragraph -------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: rework the final paragraph -------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range

# Rewrite the closing paragraph's text first. Excluding the trailing
# paragraph mark from the range being retyped keeps its italic run
# formatting (and the paragraph's leading empty run) untouched.
$bodyOnly = $d.Range($lastRange.Start, $lastRange.End - 1)
$bodyOnly.Text = "Play Fruit Awards for free and find out what's great and not so great in this classic-style slot game."

# Re-fetch the (possibly reindexed) last paragraph and split a brand new
# paragraph in front of it for the bold heading line.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$newRange = $newPara.Range
$newBodyOnly = $d.Range($newRange.Start, $newRange.End - 1)
$newBodyOnly.Text = "Play Fruit Awards | Free Review of Classic Slot Game"

$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)
$newRange = $newPara.Range
$newBodyOnly = $d.Range($newRange.Start, $newRange.End - 1)
$newBodyOnly.Font.Bold = 1
$newBodyOnly.Font.Italic = 0

Write-Output "Edit complete"
